$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2018-12-31 00:00:00"
$ws.Range("O2").Value = 29222647.67
$ws.Range("P2").Value = 13869.7235846522
$ws.Range("Q2").Value = 177551017.82
$ws.Range("R2").Value = 84269.69271046419
$ws.Range("S2").Value = 53941412.48
$ws.Range("T2").Value = 25601.8034132946
$ws.Range("U2").Value = -29846662.96
$ws.Range("V2").Value = -14165.8952280513
$ws.Range("W2").Value = 310377.55
$ws.Range("X2").Value = 147.3121420754
$ws.Range("Y2").Value = 74329537.20999999
$ws.Range("Z2").Value = 35278.4643924026
$ws.Range("AA2").Value = -7453.64
$ws.Range("AB2").Value = -3.5376646109
$ws.Range("AC2").Value = -210693.8
$ws.Range("AD2").Value = 99.147946885
